$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 9 data matching the existing table layout (A:H)
$ws.Range("A9").Value = 9895.27
$ws.Range("B9").Value = 9974.07
$ws.Range("C9").Value = 282.89999999999998
$ws.Range("D9").Value = 285.14
$ws.Range("E9").Value = $true
$ws.Range("F9").Value = 0.79

# Copy the date-formatted style from G8 down to G9, then set its value
$ws.Range("G8").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("G9").Value = 42609.487256944441

$ws.Range("H9").Value = $false
